$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking price strings from being auto-converted to numbers
$ws.Range("D2").Value = "67.347.99"
$ws.Range("E2").Value = "  -1.91%  "

$ws.Range("D3").Value = "3.705.95"
$ws.Range("E3").Value = "  -2.84%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.47"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.07"
$ws.Range("E6").Value = "  -3.10%  "

$ws.Range("D7").Value = "3.700.99"
$ws.Range("E7").Value = "  -2.77%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  +0.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").Value = "  +2.16%  "

$ws.Range("E11").Value = "  -1.81%  "

$ws.Range("E12").Value = "  -2.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.85"
$ws.Range("E13").Value = "  -2.23%  "

$ws.Range("E14").Value = "  -1.66%  "

$ws.Range("D15").Value = "4.326.36"
$ws.Range("E15").Value = "  -3.04%  "

$ws.Range("D16").Value = "3.694.44"
$ws.Range("E16").Value = "  -3.65%  "

$ws.Range("D17").Value = "67.338.97"
$ws.Range("E17").Value = "  -2.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.26"
$ws.Range("E18").Value = "  -1.34%  "

$ws.Range("E19").Value = "  -2.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.38"
$ws.Range("E20").Value = "  +7.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "487.61"
$ws.Range("E21").Value = "  -1.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.27"
$ws.Range("E22").Value = "  -1.99%  "

$ws.Range("E23").Value = "  -1.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.04"
$ws.Range("E24").Value = "  -2.07%  "

$ws.Range("E25").Value = "  +4.11%  "

$ws.Range("E26").Value = "  -3.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.21"
$ws.Range("E27").Value = "  -1.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").Value = "  -1.06%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  -1.01%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.35"
$ws.Range("E31").Value = "  -4.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.70"
$ws.Range("E32").Value = "  -0.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.41"
$ws.Range("E33").Value = "  -3.78%  "

$ws.Range("D34").Value = "3.848.85"
$ws.Range("E34").Value = "  -3.07%  "

$ws.Range("E35").Value = "  -3.30%  "

$ws.Range("D36").Value = "3.651.00"
$ws.Range("E36").Value = "  -2.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -2.31%  "

$ws.Range("E39").Value = "  -1.91%  "

$ws.Range("E40").Value = "  -3.78%  "

$ws.Range("E41").Value = "  -1.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "48.74"
$ws.Range("E42").Value = "  -1.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "426.82"
$ws.Range("E43").Value = "  -6.06%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.93"
$ws.Range("E44").Value = "  -4.41%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.83"
$ws.Range("E45").Value = "  -1.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.43"
$ws.Range("E46").Value = "  +0.00%  "

$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.46"
$ws.Range("E48").Value = "  -2.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.73"
$ws.Range("E49").Value = "  +1.44%  "

$ws.Range("D50").Value = "2.754.08"
$ws.Range("E50").Value = "  -3.81%  "

$ws.Range("E51").Value = "  -1.72%  "

